# "format dan tugas bab 2"
#
# 1) Remove the stray _GoBack bookmark that currently sits right after
#    "...disertai dengan penjelasan masing masing kelompok." (end of a
#    NormalWeb paragraph).
# 2) Fix the deadline date ("Deadline ACC 9 November 2018" -> "...6...")
#    and leave a fresh _GoBack bookmark immediately after the corrected
#    digit, matching where Word drops _GoBack after the last edit.

$d = $word.ActiveDocument

# --- 1. Drop the old _GoBack bookmark -------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Locate the standalone "9" in "Deadline ACC 9 November 2018" -------
$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$found = $rng.Find.Execute("9", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $digitStart = $rng.Start
    $digitEnd = $rng.End

    # Temporary barrier bookmark right before the digit so rewriting its
    # text does not get coalesced into the preceding " " run.
    $barrierRange = $d.Range($digitStart, $digitStart)
    $d.Bookmarks.Add("ZZZTempBarrier", $barrierRange)

    # The real _GoBack bookmark goes right after the digit, collapsed
    # (i.e. it does not wrap the run), which also stops it merging with
    # the following " November 2018" run.
    $goBackRange = $d.Range($digitEnd, $digitEnd)
    $d.Bookmarks.Add("_GoBack", $goBackRange)

    # Now change the digit itself: 9 -> 6.
    $digitRange = $d.Range($digitStart, $digitEnd)
    $digitRange.Text = "6"

    # Clean up the temporary barrier bookmark.
    $d.Bookmarks("ZZZTempBarrier").Delete()
}
